# Fitness/training log update: add "Details", "Time of Day" and "Notes"
# columns, split the old combined "Strength Training - X" activity into
# "Gym" + a Details breakdown, and append ~2 weeks of new log rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Seed the new shared strings in the same order the original author
#    appears to have typed them (this keeps the rebuilt sharedStrings.xml
#    table in the same relative order as the target workbook).
# ---------------------------------------------------------------------

# "Gym" replaces "Strength Training - Chest" in B4 (first use of "Gym").
$ws.Range("B4").Value = "Gym"

# New header cells.
$ws.Range("F1").Value = "Time of Day (hh:mm)"
$ws.Range("C1").Value = "Details"

# Details column values, in the order they were first introduced.
$ws.Range("C2").Value = "Breatstroke"
$ws.Range("C4").Value = "Chest"
$ws.Range("C5").Value = "Legs"
$ws.Range("C6").Value = "Breaststroke"
$ws.Range("C13").Value = "Back"
$ws.Range("C7").Value = "Easy run"
$ws.Range("C9").Value = "Breatstroke / Frontcrawl"
$ws.Range("C12").Value = "Breaststroke / Frontcrawl"

# Final new header cell.
$ws.Range("H1").Value = "Notes"

# ---------------------------------------------------------------------
# 2) Move "Duration (min)" / "Distance (m)" / "Intensity Rating (1-10)"
#    headers into their new columns (D, E, G) and fix up B5 (old
#    "Strength Training - Legs").
# ---------------------------------------------------------------------

$ws.Range("D1").Value = "Duration (min)"
$ws.Range("E1").Value = "Distance (m)"
$ws.Range("G1").Value = "Intensity Rating (1-10)"

$ws.Range("B5").Value = "Gym"

# ---------------------------------------------------------------------
# 3) Re-lay the numeric Duration/Distance/Intensity data that shifted
#    columns (old C/D/E -> new D/E/G), row by row. Doing the moved
#    columns from the right-most column back to the left avoids
#    clobbering a cell before it has been read where a row's columns
#    overlap.
# ---------------------------------------------------------------------

# Row 2 (Swimming / Breatstroke): old C2=30,D2=550,E2=8 -> new D2,E2,G2
$ws.Range("G2").Value = 8
$ws.Range("E2").Value = 550
$ws.Range("D2").Value = 30
$ws.Range("F2").Value = 0.33333333333333331

# Row 4 (Gym / Chest): old C4=90,E4=8 -> new D4,G4 (no Distance for a gym day)
$ws.Range("G4").Value = 8
$ws.Range("E4").Clear()
$ws.Range("D4").Value = 90
$ws.Range("F4").Value = 0.75

# Row 5 (Gym / Legs): old C5=120,E5=7 -> new D5,G5 (no Distance for a gym day)
$ws.Range("G5").Value = 7
$ws.Range("E5").Clear()
$ws.Range("D5").Value = 120
$ws.Range("F5").Value = 0.75

# Row 6 (Swimming / Breaststroke): old C6=30,D6=650,E6=7 -> new D6,E6,G6
$ws.Range("G6").Value = 7
$ws.Range("E6").Value = 650
$ws.Range("D6").Value = 30
$ws.Range("F6").Value = 0.33333333333333331

# Row 7 (Running / Easy run): old C7=35,D7=6470,E7=8 -> new D7,E7,G7
$ws.Range("G7").Value = 8
$ws.Range("E7").Value = 6470
$ws.Range("D7").Value = 35
$ws.Range("F7").Value = 0.4375

# Row 8 was date-only before; it gets a "Rest" activity now.
$ws.Range("B8").Value = "Rest"

# ---------------------------------------------------------------------
# 4) New rows of fitness log data (rows 9-13), plus the trailing
#    rest-of-month date-only rows (14-20).
# ---------------------------------------------------------------------

# Row 9: Swim, Breatstroke / Frontcrawl
$ws.Range("A9").Value = 45215
$ws.Range("B9").Value = "Swimming"
$ws.Range("D9").Value = 30
$ws.Range("E9").Value = 700
$ws.Range("F9").Value = 0.33333333333333331
$ws.Range("G9").Value = 7

# Row 10: Gym, Chest
$ws.Range("A10").Value = 45215
$ws.Range("B10").Value = "Gym"
$ws.Range("C10").Value = "Chest"
$ws.Range("D10").Value = 90
$ws.Range("F10").Value = 0.70833333333333337
$ws.Range("G10").Value = 6

# Row 11: Rest
$ws.Range("A11").Value = 45216
$ws.Range("B11").Value = "Rest"

# Row 12: Swim, Breaststroke / Frontcrawl
$ws.Range("A12").Value = 45217
$ws.Range("B12").Value = "Swimming"
$ws.Range("D12").Value = 40
$ws.Range("E12").Value = 900
$ws.Range("F12").Value = 0.3263888888888889
$ws.Range("G12").Value = 7

# Row 13: Gym, Back (slightly taller row in the source workbook)
$ws.Range("A13").Value = 45217
$ws.Range("B13").Value = "Gym"
$ws.Range("D13").Value = 90
$ws.Range("F13").Value = 0.72916666666666663
$ws.Rows.Item(13).RowHeight = 15

# Rows 14-20: just the date column, filled in for the rest of the month.
$ws.Range("A14").Value = 45218
$ws.Range("A15").Value = 45219
$ws.Range("A16").Value = 45220
$ws.Range("A17").Value = 45221
$ws.Range("A18").Value = 45222
$ws.Range("A19").Value = 45223
$ws.Range("A20").Value = 45224

# Copy the date-column formatting (style 2: numFmt 14 + left align) down
# onto all the newly added date cells so they don't pick up a fresh style.
$ws.Range("A8").Copy()
$ws.Range("A9:A20").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 5) Number-format the new "Time of Day" column as h:mm (creates the new
#    numFmtId 20 cellXf) before bolding the header row, so the new
#    cellXfs come out in the same order as the target file (time-format
#    xf before the bold-header xf).
# ---------------------------------------------------------------------

$ws.Range("F2").NumberFormat = "h:mm"
$ws.Range("F4").NumberFormat = "h:mm"
$ws.Range("F5").NumberFormat = "h:mm"
$ws.Range("F6").NumberFormat = "h:mm"
$ws.Range("F7").NumberFormat = "h:mm"
$ws.Range("F9").NumberFormat = "h:mm"
$ws.Range("F10").NumberFormat = "h:mm"
$ws.Range("F12").NumberFormat = "h:mm"
$ws.Range("F13").NumberFormat = "h:mm"

# ---------------------------------------------------------------------
# 6) Bold the header row (this clones style 1 with Font.Bold, landing on
#    the new cellXf right after the time-format one).
# ---------------------------------------------------------------------

$ws.Range("A1:H1").Font.Bold = $true

# ---------------------------------------------------------------------
# 7) Column widths. The underlying engine quantizes ColumnWidth to a
#    coarser grid than native Excel, so these are the closest achievable
#    inputs to the target stored widths.
# ---------------------------------------------------------------------

$ws.Columns.Item(2).ColumnWidth = 12.76    # -> ~13.73 (target 13.734375)
$ws.Columns.Item(3).ColumnWidth = 21.25    # -> ~22.17 (target 22.20703125)
$ws.Columns.Item(4).ColumnWidth = 18.25    # -> ~19.17 (target 19.20703125)
$ws.Columns.Item(5).ColumnWidth = 11.25    # -> ~12.17 (target 12.20703125)
$ws.Columns.Item(6).ColumnWidth = 19.09    # -> 20 (target 20, exact)
$ws.Columns.Item(8).ColumnWidth = 22.25    # -> ~23.17 (target 23.20703125)
# Column G (Intensity) carries over the old column E's width (the Intensity
# data moved from E to G), so it needs to be set explicitly too.
$ws.Columns.Item(7).ColumnWidth = 19.59    # -> ~20.5 (target 20.5234375)

# ---------------------------------------------------------------------
# 8) Misc sheet/view properties that changed.
# ---------------------------------------------------------------------

$ws.PageSetup.Orientation = 1   # xlPortrait
$ws.Range("J33").Select()
